$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT (matches source data which stores
# numeric-looking price strings, e.g. "60.893.26", as literal text) without
# leaving a stray number-format style behind on the cell.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value2 = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "60.893.26"
Set-TextValue "E2" "  +0.28%  "

# Row 3
Set-TextValue "D3" "2.594.13"
Set-TextValue "E3" "  +0.15%  "

# Row 4
Set-TextValue "E4" "  +0.09%  "

# Row 5
Set-TextValue "D5" "523.35"
Set-TextValue "E5" "  +2.98%  "

# Row 6
Set-TextValue "D6" "154.37"
Set-TextValue "E6" "  +0.20%  "

# Row 7
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  +0.10%  "

# Row 8
Set-TextValue "D8" "0.593"
Set-TextValue "E8" "  +2.20%  "

# Row 9
Set-TextValue "D9" "6.70"
Set-TextValue "E9" "  +1.66%  "

# Row 10
Set-TextValue "E10" "  +1.60%  "

# Row 11
Set-TextValue "E11" "  +0.06%  "

# Row 12
Set-TextValue "D12" "0.129"
Set-TextValue "E12" "  +1.40%  "

# Row 13
Set-TextValue "D13" "3.050.81"
Set-TextValue "E13" "  +0.33%  "

# Row 14
Set-TextValue "D14" "60.938.43"
Set-TextValue "E14" "  +0.56%  "

# Row 15
Set-TextValue "D15" "21.64"
Set-TextValue "E15" "  +0.16%  "

# Row 16
Set-TextValue "E16" "  +0.03%  "

# Row 17
Set-TextValue "D17" "2.602.23"
Set-TextValue "E17" "  +0.53%  "

# Row 19
Set-TextValue "D19" "353.75"
Set-TextValue "E19" "  +2.11%  "

# Row 20
Set-TextValue "D20" "10.59"
Set-TextValue "E20" "  +1.26%  "

# Row 21
Set-TextValue "D21" "6.24"
Set-TextValue "E21" "  +1.67%  "

# Row 22
Set-TextValue "D22" "0.999"
Set-TextValue "E22" "  +0.21%  "

# Row 23
Set-TextValue "D23" "60.94"
Set-TextValue "E23" "  +1.46%  "

# Row 24
Set-TextValue "E24" "  +1.55%  "

# Row 25
Set-TextValue "E25" "  -0.93%  "

# Row 26
Set-TextValue "D26" "2.710.55"
Set-TextValue "E26" "  +0.34%  "

# Row 27
Set-TextValue "D27" "0.997"
Set-TextValue "E27" "  -0.03%  "

# Row 28
Set-TextValue "E28" "  -0.45%  "

# Row 29
Set-TextValue "E29" "  +0.04%  "

# Row 31
Set-TextValue "D31" "6.34"
Set-TextValue "E31" "  +10.95%  "

# Row 32
Set-TextValue "D32" "19.37"
Set-TextValue "E32" "  +0.04%  "

# Row 33
Set-TextValue "E33" "  +2.80%  "

# Row 34
Set-TextValue "D34" "148.29"
Set-TextValue "E34" "  -3.43%  "

# Row 35
Set-TextValue "D35" "4.16"
Set-TextValue "E35" "  +4.10%  "

# Row 36
Set-TextValue "D36" "0.939"
Set-TextValue "E36" "  +8.89%  "

# Row 37
Set-TextValue "E37" "  +0.65%  "

# Row 38
Set-TextValue "E38" "  +1.77%  "

# Row 39
Set-TextValue "B39" "Filecoin"
Set-TextValue "C39" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D39" "3.80"
Set-TextValue "E39" "  +0.84%  "

# Row 40
Set-TextValue "B40" "Fetch.AI"
Set-TextValue "C40" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D40" "0.849"
Set-TextValue "E40" "  -0.65%  "

# Row 41
Set-TextValue "B41" "OKB"
Set-TextValue "C41" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D41" "36.45"
Set-TextValue "E41" "  +1.44%  "

# Row 42
Set-TextValue "D42" "288.47"
Set-TextValue "E42" "  -2.80%  "

# Row 43
Set-TextValue "E43" "  +1.55%  "

# Row 44
Set-TextValue "E44" "  +1.58%  "

# Row 45
Set-TextValue "E45" "  +0.15%  "

# Row 46
Set-TextValue "D46" "0.998"
Set-TextValue "E46" "  +0.18%  "

# Row 47
Set-TextValue "D47" "19.61"
Set-TextValue "E47" "  -1.22%  "

# Row 48
Set-TextValue "B48" "VeChain"
Set-TextValue "C48" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D48" "0.0238"
Set-TextValue "E48" "  +2.09%  "

# Row 49
Set-TextValue "B49" "RenderToken"
Set-TextValue "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "4.88"
Set-TextValue "E49" "  +0.17%  "

# Row 50
Set-TextValue "D50" "10.33"
Set-TextValue "E50" "  +0.19%  "

# Row 51
Set-TextValue "D51" "19.06"
Set-TextValue "E51" "  +7.84%  "

